# Add two new rows (17 and 18) to the LeetCode pandas pattern tracker sheet,
# and expand the Table2 ListObject + autofilter to cover them.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 17: 1050. Actors and Directors Who Cooperated At Least Three Times ---

# Seed row 17 with the same formatting as row 16 (the last existing data row)
# so fills/fonts carry over correctly before we fill in the new values.
$ws.Range("A16:E16").Copy()
$ws.Range("A17:E17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A17").Value = "1050. Actors and Directors Who Cooperated At Least Three Times"
$ws.Range("B17").Value = "Easy"
$ws.Range("C17").Value = "Data Integration"
$ws.Range("D17").Value = "Use groupyby(actor, director).agg(count=(director, count)), reset index, then return the stats df stats[stats['count']>=3][['actor_id]', 'director_id']"
$ws.Range("E17").Value = "https://leetcode.com/problems/actors-and-directors-who-cooperated-at-least-three-times/solutions/3946925/pandas-2-step-simple-code-with-comments/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "

$ws.Hyperlinks.Add($ws.Range("E17"), "https://leetcode.com/problems/actors-and-directors-who-cooperated-at-least-three-times/solutions/3946925/pandas-2-step-simple-code-with-comments/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata")

# Restore the normal hyperlink-column formatting (Hyperlinks.Add re-styles the cell).
$ws.Range("E16").Copy()
$ws.Range("E17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Row 18: 607. Sales Person ---

$ws.Range("A17:E17").Copy()
$ws.Range("A18:E18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("A18").Value = "607. Sales Person"
$ws.Range("B18").Value = "Easy"
$ws.Range("C18").Value = "Data Integration"
$ws.Range("D18").Value = "merge the sales_person data with orders on sales id, how=left, and merge company on com_id left, filter data where name_y is not RED or null (.isna()), rename name_x column to name in result df, drop duplicates, then find the salespersons who had orders related to RED and name_x with .unique(), then filter out the salespersons with RED orders from the result."
$ws.Range("E18").Value = "https://leetcode.com/problems/sales-person/solutions/3970324/pandas-easy-to-understand-properly-commented/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata "

$ws.Hyperlinks.Add($ws.Range("E18"), "https://leetcode.com/problems/sales-person/solutions/3970324/pandas-easy-to-understand-properly-commented/?envType=study-plan-v2&envId=30-days-of-pandas&lang=pythondata")

$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Expand the table (Table2) and its autofilter to include the new rows ---
$tbl = $ws.ListObjects.Item(1)
$tbl.Resize($ws.Range("A1:E18"))

# --- Update the active cell selection to match the saved workbook state ---
$ws.Range("E25").Select()

Write-Output "Edit complete"
